$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row pair (4,5): swap index/prolificid/name for the two workers, and
# refresh the recomputed re_rank (realeffort) scores for every Asian row.
$ws.Range("G2").Value = 11.36491441729315
$ws.Range("G3").Value = 10.01785415257338

$ws.Range("C4").Value = 2
$ws.Range("D4").Value = "5f2c1a97a6809c060fec8820"
$ws.Range("E4").Value = "Maggie"
$ws.Range("G4").Value = 8.040674606944371

$ws.Range("C5").Value = 10
$ws.Range("D5").Value = "60a71d27a66fac796ad4de6f"
$ws.Range("E5").Value = "Jennifer"
$ws.Range("G5").Value = 8.032360915298707

$ws.Range("G6").Value = 7.393130267324382
$ws.Range("G7").Value = 6.317503956260554
$ws.Range("G8").Value = 6.011467763540303
$ws.Range("G9").Value = 5.079227112452898
$ws.Range("G10").Value = 5.054581671041178
$ws.Range("G11").Value = 4.380943723260605
$ws.Range("G12").Value = 2.380573180982641
$ws.Range("G13").Value = 0.2369653110842641

# Row pair (14,15): swap index/prolificid/name.
$ws.Range("C14").Value = 7
$ws.Range("D14").Value = "6024c18b094ac71dd93f4f5a"
$ws.Range("E14").Value = "Katherine"
$ws.Range("G14").Value = 8.051697533201137

$ws.Range("C15").Value = 2
$ws.Range("D15").Value = "60778ed0fde3e9c3a96f1d11"
$ws.Range("E15").Value = "Melissa"
$ws.Range("G15").Value = 8.049345038247747

# Row pair (16,17): swap index/prolificid/name/gender.
$ws.Range("C16").Value = 8
$ws.Range("D16").Value = "5f0142aa1eb1e528e7abce50"
$ws.Range("E16").Value = "Valeria"
$ws.Range("F16").Value = "female"
$ws.Range("G16").Value = 7.172380869265427

$ws.Range("C17").Value = 3
$ws.Range("D17").Value = "60ba8ba51a5e0a105396888a"
$ws.Range("E17").Value = "Alfredo"
$ws.Range("F17").Value = "male"
$ws.Range("G17").Value = 7.084258182079134

# Row pair (18,19): swap index/prolificid/name.
$ws.Range("C18").Value = 11
$ws.Range("D18").Value = "5f5ea8227fa75676f56f9276"
$ws.Range("E18").Value = "Carlos"
$ws.Range("G18").Value = 6.251133082574972

$ws.Range("C19").Value = 0
$ws.Range("D19").Value = "5eeaa065c7acf61c4322f6d9"
$ws.Range("E19").Value = "Yonifredy"
$ws.Range("G19").Value = 6.228542514609791

$ws.Range("G20").Value = 5.257951082805501

# Row pair (21,22): swap index/prolificid/name.
$ws.Range("C21").Value = 4
$ws.Range("D21").Value = "5e706891c396cc64388ef760"
$ws.Range("E21").Value = "Maria"
$ws.Range("G21").Value = 3.374220793453376

$ws.Range("C22").Value = 1
$ws.Range("D22").Value = "5e0adc8f4cac6834756db412"
$ws.Range("E22").Value = "Mary"
$ws.Range("G22").Value = 3.30059129783378

$ws.Range("G23").Value = 2.445237386649406
$ws.Range("G24").Value = 1.109097411672669
$ws.Range("G25").Value = 0.3238469706846104
